$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 246 (old rows 246-253 shift down to 249-256)
$ws.Rows("246:248").Insert()

# Row 246 (new)
$ws.Cells.Item(246,1).Value = 4
$ws.Cells.Item(246,2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(246,3).Value = "Los Lagos"
$ws.Cells.Item(246,4).Value = 44610
$ws.Cells.Item(246,5).Value = 10
$ws.Cells.Item(246,6).Value = 100112028
$ws.Cells.Item(246,7).Value = "Sandia"
$ws.Cells.Item(246,8).Value = "Sin especificar"
$ws.Cells.Item(246,9).Value = "Primera"
$ws.Cells.Item(246,10).Value = 1000
$ws.Cells.Item(246,11).Value = 3000
$ws.Cells.Item(246,12).Value = 3000
$ws.Cells.Item(246,13).Value = 3000
$ws.Cells.Item(246,14).Value = "`$/unidad"
$ws.Cells.Item(246,15).Value = "Región de O'Higgins"
$ws.Cells.Item(246,16).Value = 3000
$ws.Cells.Item(246,17).Value = 1
$ws.Cells.Item(246,18).Value = "Hortaliza"

# Row 247 (new)
$ws.Cells.Item(247,1).Value = 4
$ws.Cells.Item(247,2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(247,3).Value = "Los Lagos"
$ws.Cells.Item(247,4).Value = 44610
$ws.Cells.Item(247,5).Value = 10
$ws.Cells.Item(247,6).Value = 100112028
$ws.Cells.Item(247,7).Value = "Sandia"
$ws.Cells.Item(247,8).Value = "Sin especificar"
$ws.Cells.Item(247,9).Value = "Segunda"
$ws.Cells.Item(247,10).Value = 1000
$ws.Cells.Item(247,11).Value = 2500
$ws.Cells.Item(247,12).Value = 2500
$ws.Cells.Item(247,13).Value = 2500
$ws.Cells.Item(247,14).Value = "`$/unidad"
$ws.Cells.Item(247,15).Value = "Región de O'Higgins"
$ws.Cells.Item(247,16).Value = 2500
$ws.Cells.Item(247,17).Value = 1
$ws.Cells.Item(247,18).Value = "Hortaliza"

# Row 248 (new)
$ws.Cells.Item(248,1).Value = 4
$ws.Cells.Item(248,2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(248,3).Value = "Los Lagos"
$ws.Cells.Item(248,4).Value = 44610
$ws.Cells.Item(248,5).Value = 10
$ws.Cells.Item(248,6).Value = 100112028
$ws.Cells.Item(248,7).Value = "Sandia"
$ws.Cells.Item(248,8).Value = "Sin especificar"
$ws.Cells.Item(248,9).Value = "Tercera"
$ws.Cells.Item(248,10).Value = 1500
$ws.Cells.Item(248,11).Value = 2000
$ws.Cells.Item(248,12).Value = 2000
$ws.Cells.Item(248,13).Value = 2000
$ws.Cells.Item(248,14).Value = "`$/unidad"
$ws.Cells.Item(248,15).Value = "Región de O'Higgins"
$ws.Cells.Item(248,16).Value = 2000
$ws.Cells.Item(248,17).Value = 1
$ws.Cells.Item(248,18).Value = "Hortaliza"
